$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 11202
$ws.Range("E2").Value = -165
$ws.Range("F2").Value = -166
$ws.Range("G2").Value = -2606
$ws.Range("H2").Value = -4127
$ws.Range("I2").Value = -4120
$ws.Range("J2").Value = -7
$ws.Range("K2").Value = 18708
$ws.Range("L2").Value = 18562
$ws.Range("M2").Value = 146
$ws.Range("N2").Value = 113
$ws.Range("O2").Value = 34
$ws.Range("P2").Value = 204
$ws.Range("Q2").Value = -356
$ws.Range("R2").Value = 461
$ws.Range("S2").Value = -343
$ws.Range("T2").Value = 54
$ws.Range("U2").Value = -410
$ws.Range("V2").Value = 15187
$ws.Range("W2").Value = -1.48
$ws.Range("X2").Value = -36.84
$ws.Range("Y2").Value = -183.3
$ws.Range("Z2").Value = -19.94
$ws.Range("AA2").Value = 12674.97
$ws.Range("AB2").Value = -346.8
$ws.Range("AC2").Value = -100795
$ws.Range("AD2").Value = -0.08
$ws.Range("AE2").Value = 2762
$ws.Range("AF2").Value = 2.8
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 4087955

# Row 3
$ws.Range("D3").Value = 8527
$ws.Range("E3").Value = -415
$ws.Range("F3").Value = -349
$ws.Range("G3").Value = -1071
$ws.Range("H3").Value = -1260
$ws.Range("I3").Value = -1259
$ws.Range("J3").Value = -1
$ws.Range("K3").Value = 10763
$ws.Range("L3").Value = 9597
$ws.Range("M3").Value = 1166
$ws.Range("N3").Value = 1136
$ws.Range("O3").Value = 30
$ws.Range("P3").Value = 2124
$ws.Range("Q3").Value = -354
$ws.Range("R3").Value = 5831
$ws.Range("S3").Value = -5828
$ws.Range("T3").Value = 77
$ws.Range("U3").Value = -431
$ws.Range("V3").Value = 6850
$ws.Range("W3").Value = -4.87
$ws.Range("X3").Value = -14.78
$ws.Range("Y3").Value = -201.64
$ws.Range("Z3").Value = -8.550000000000001
$ws.Range("AA3").Value = 822.98
$ws.Range("AB3").Value = -73.02
$ws.Range("AC3").Value = -3483
$ws.Range("AD3").Value = -1
$ws.Range("AE3").Value = 2775
$ws.Range("AF3").Value = 1.26
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 24087955

# Row 4
$ws.Range("D4").Value = 8370
$ws.Range("E4").Value = -212
$ws.Range("F4").Value = -201
$ws.Range("G4").Value = -1519
$ws.Range("H4").Value = -1349
$ws.Range("I4").Value = -1007
$ws.Range("J4").Value = -342
$ws.Range("K4").Value = 11086
$ws.Range("L4").Value = 9084
$ws.Range("M4").Value = 2002
$ws.Range("N4").Value = 2012
$ws.Range("O4").Value = -10
$ws.Range("P4").Value = 90
$ws.Range("Q4").Value = 387
$ws.Range("R4").Value = 469
$ws.Range("S4").Value = -823
$ws.Range("T4").Value = 55
$ws.Range("U4").Value = 332
$ws.Range("V4").Value = 5950
$ws.Range("W4").Value = -2.54
$ws.Range("X4").Value = -16.12
$ws.Range("Y4").Value = -63.97
$ws.Range("Z4").Value = -12.35
$ws.Range("AA4").Value = 453.72
$ws.Range("AB4").Value = 5127.05
$ws.Range("AC4").Value = -2363
$ws.Range("AD4").Value = -1.84
$ws.Range("AE4").Value = 3950
$ws.Range("AF4").Value = 1.1
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 24087955

# Row 5
$ws.Range("D5").Value = 9003
$ws.Range("E5").Value = -120
$ws.Range("F5").Value = -120
$ws.Range("G5").Value = 2042
$ws.Range("H5").Value = 2231
$ws.Range("I5").Value = 2306
$ws.Range("J5").Value = -74
$ws.Range("K5").Value = 11141
$ws.Range("L5").Value = 7105
$ws.Range("M5").Value = 4036
$ws.Range("N5").Value = 4303
$ws.Range("O5").Value = -267
$ws.Range("P5").Value = 2262
$ws.Range("Q5").Value = 102
$ws.Range("R5").Value = 3655
$ws.Range("S5").Value = -3639
$ws.Range("T5").Value = 41
$ws.Range("U5").Value = 61
$ws.Range("V5").Value = 3316
$ws.Range("W5").Value = -1.34
$ws.Range("X5").Value = 24.79
$ws.Range("Y5").Value = 73.02
$ws.Range("Z5").Value = 20.08
$ws.Range("AA5").Value = 176.03
$ws.Range("AB5").Value = 505.56
$ws.Range("AC5").Value = 7009
$ws.Range("AD5").Value = 0.55
$ws.Range("AE5").Value = 12955
$ws.Range("AF5").ClearContents()
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").Value = 0.3
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 45235478

# Row 6
$ws.Range("D6").Value = 9455
$ws.Range("E6").Value = 61
$ws.Range("F6").Value = 61
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = -68
$ws.Range("I6").Value = -90
$ws.Range("K6").Value = 10928
$ws.Range("L6").Value = 6414
$ws.Range("M6").Value = 4514
$ws.Range("N6").Value = 4762
$ws.Range("P6").Value = 2262
$ws.Range("Q6").Value = -509
$ws.Range("R6").Value = 710
$ws.Range("S6").Value = -112
$ws.Range("T6").Value = 69
$ws.Range("U6").Value = -578
$ws.Range("V6").Value = 3593
$ws.Range("W6").Value = 0.65
$ws.Range("X6").Value = -0.72
$ws.Range("Y6").Value = -1.98
$ws.Range("Z6").Value = -0.62
$ws.Range("AA6").Value = 142.1
$ws.Range("AB6").Value = 518.5
$ws.Range("AC6").Value = -198
$ws.Range("AD6").Value = -20.95
$ws.Range("AE6").Value = 14335
$ws.Range("AF6").ClearContents()
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").Value = 0.29
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 45235478

# Row 7 - clear all data cells except A,B,C
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8 - clear all data cells except A,B,C
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9 - clear all data cells except A,B,C
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
